# "arreglo de letras cortadas" - fix truncated abbreviation letters in the
# "Extension" (H) column: these academic-degree abbreviations were missing
# their trailing period, which made them look cut off. Restore it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tijuana")

$ws.Range("H2").Value  = "L.A.E."
$ws.Range("H4").Value  = "L.M.I."
$ws.Range("H6").Value  = "L.C.E.A."
$ws.Range("H7").Value  = "L.C.E.A."
$ws.Range("H8").Value  = "L.C."
$ws.Range("H12").Value = "L.C."
$ws.Range("H14").Value = "L.C.E.A."
$ws.Range("H18").Value = "I.D.I.E."
$ws.Range("H19").Value = "L.C.E.A."

# View-state touch-ups: the author re-scrolled/re-zoomed and left the
# selection on H20 before saving.
$ws.Activate()
$ws.Range("H20").Select()
$excel.ActiveWindow.Zoom = 100
